# Generate Report for Handoff
# Update status + handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (B2, C2) move from "In Translation" to "Ready for handoff",
# and the Latest Handoff Date (D2) is refreshed.
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-29-18 14:29:23"

# zh-cn sheet: Latest Handoff Datetime (E2) refreshed.
$zhcn.Range("E2").Value = "2016-03-18 14:29:20"

# de-de sheet: Latest Handoff Datetime (E2) refreshed.
$dede.Range("E2").Value = "2016-03-18 14:29:23"
